$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -12.1312
$ws.Range("B7").Value = 4.7643
$ws.Range("E7").Value = 15.93390000000001
$ws.Range("A9").Value = -22.0013
$ws.Range("E10").Value = 16.82709999999999
$ws.Range("B12").Value = 5.504499999999996
$ws.Range("A13").Value = -22.471
$ws.Range("E13").Value = 16.55130000000001
$ws.Range("B14").Value = 5.641800000000002
$ws.Range("C15").Value = -13.8136
$ws.Range("A16").Value = -21.66109999999999
$ws.Range("E16").Value = 16.19790000000001
$ws.Range("A18").Value = -22.37580000000003
$ws.Range("B19").Value = 8.655900000000004
$ws.Range("A20").Value = -19.60429999999999
$ws.Range("E20").Value = 16.29939999999999
$ws.Range("E24").Value = 16.48080000000001
$ws.Range("A26").Value = -21.16979999999997
$ws.Range("B26").Value = 4.341900000000002
$ws.Range("A27").Value = -21.47019999999996
$ws.Range("B27").Value = 5.484800000000003
$ws.Range("C28").Value = -13.4734
$ws.Range("A29").Value = -21.69819999999998
$ws.Range("B29").Value = 5.264599999999998
$ws.Range("E32").Value = 16.07419999999998
$ws.Range("C33").Value = -11.37589999999999
$ws.Range("A35").Value = -19.9293
$ws.Range("C35").Value = -12.2124
$ws.Range("A36").Value = -19.8576
$ws.Range("B37").Value = 9.1873
$ws.Range("B38").Value = 4.941800000000002
$ws.Range("C38").Value = -11.77910000000001
$ws.Range("E39").Value = 16.2824
$ws.Range("C43").Value = -14.45459999999999
$ws.Range("C44").Value = -13.86699999999999
$ws.Range("A45").Value = -21.94779999999999
$ws.Range("C45").Value = -13.82879999999999
$ws.Range("B47").Value = 5.431500000000001
$ws.Range("C47").Value = -12.8895
$ws.Range("E47").Value = 16.78199999999999
$ws.Range("E48").Value = 17.56749999999999
$ws.Range("B51").Value = 5.826600000000004
$ws.Range("C51").Value = -11.84549999999999
$ws.Range("B52").Value = 5.435699999999996
$ws.Range("E52").Value = 17.0816
$ws.Range("C54").Value = -13.15879999999999
$ws.Range("A55").Value = -22.51940000000001
$ws.Range("B55").Value = 4.727299999999996
$ws.Range("E56").Value = 16.7863
$ws.Range("A57").Value = -21.84269999999999
$ws.Range("C57").Value = -13.26559999999999
$ws.Range("C62").Value = -14.11000000000001
$ws.Range("C63").Value = -11.088
$ws.Range("C67").Value = -10.4282
$ws.Range("A69").Value = -21.66869999999999
$ws.Range("B69").Value = 5.546499999999996
$ws.Range("B70").Value = 5.868500000000004
$ws.Range("C70").Value = -11.47279999999999
$ws.Range("A76").Value = -22.34030000000001
$ws.Range("B76").Value = 5.326299999999995
$ws.Range("A78").Value = -19.99299999999997
$ws.Range("B81").Value = 5.290100000000001
$ws.Range("C81").Value = -12.12749999999999
$ws.Range("A82").Value = -21.97000000000001
$ws.Range("A83").Value = -21.8916
$ws.Range("B83").Value = 6.189900000000005
$ws.Range("E84").Value = 17.13319999999998
$ws.Range("C88").Value = -12.36269999999999
$ws.Range("A93").Value = -20.47789999999997
$ws.Range("B94").Value = 5.100999999999998
$ws.Range("C96").Value = -12.2347
$ws.Range("A97").Value = -21.69240000000001
$ws.Range("C99").Value = -11.9482
$ws.Range("B100").Value = 5.981500000000003
$ws.Range("E100").Value = 16.3958
$ws.Range("E101").Value = 16.90840000000001
$ws.Range("B102").Value = 9.234400000000004
